# Apply the diff: strip stray internal spaces (and full-width commas) that
# were accidentally introduced into several shared strings across the
# workbook's sheets. Only the textual content of the cells changes; no
# structural changes are required.

$wb = $excel.ActiveWorkbook

# --- Sheet "汽車" (car) ---
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("E2").Value = "94年05月03日"

# --- Sheet "存款" (deposits) ---
$wsDeposit = $wb.Worksheets.Item("存款")
$wsDeposit.Range("B3").Value = "匯豐（台灣）商業銀行台北分行"
$wsDeposit.Range("B4").Value = "台新國際商業銀行板橋分行"
$wsDeposit.Range("B5").Value = "遠東國際商業銀行板橋埔墘分行"
$wsDeposit.Range("B6").Value = "台北台大郵局（第23支局）"
$wsDeposit.Range("B7").Value = "台北台大郵局（第23支局）"
$wsDeposit.Range("B8").Value = "遠東國際商業銀行板橋埔墘分行"

# --- Sheet "債務" (debt) ---
$wsDebt = $wb.Worksheets.Item("債務")
$wsDebt.Range("D2").Value = "遠東銀行板橋中正分行新北市板橋區中正路228號"
$wsDebt.Range("F2").Value = "102年08月20日"
$wsDebt.Range("D3").Value = "遠東銀行板橋中正分行新北市板橋區中正路228號"
# E3 held a text value ("3，754，458") that merely had its full-width commas
# removed; it must stay text-typed (not become numeric), so force the
# number format to Text before writing the digits-only string.
$wsDebt.Range("E3").NumberFormat = "@"
$wsDebt.Range("E3").Value = "3754458"
$wsDebt.Range("F3").Value = "102年08月20日"
